$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("M31").Value = " "
Write-Host $ws.Range("M31").Value
